$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 482.94116
$ws.Range("J17").Value = 482.94116
$ws.Range("L17").Value = 1448.82348
$ws.Range("N17").Value = -1784.82348
$ws.Range("H43").Value = 3311.0667
$ws.Range("J43").Value = 3462.125
$ws.Range("L43").Value = 3462.125
$ws.Range("N43").Value = -3600.125
$ws.Range("H64").Value = 4855.857
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("H67").Value = 4855.857
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("H74").Value = 38230.418
$ws.Range("I74").Value = 49474.223
$ws.Range("K74").Value = 49474.223
$ws.Range("M74").Value = -48538.223
$ws.Range("H76").Value = 3976.6365
$ws.Range("I76").Value = 3737.5
$ws.Range("J76").Value = 4113.2856
$ws.Range("K76").Value = 3737.5
$ws.Range("L76").Value = 4113.2856
$ws.Range("M76").Value = -3422.5
$ws.Range("N76").Value = -4743.2856
$ws.Range("H77").Value = 38230.418
$ws.Range("I77").Value = 49474.223
$ws.Range("K77").Value = 247371.115
$ws.Range("M77").Value = -242691.115
$ws.Range("H79").Value = 3976.6365
$ws.Range("I79").Value = 3737.5
$ws.Range("J79").Value = 4113.2856
$ws.Range("K79").Value = 3737.5
$ws.Range("L79").Value = 4113.2856
$ws.Range("M79").Value = -2645.5
$ws.Range("N79").Value = -6297.2856
$ws.Range("H116").Value = 6603.091
$ws.Range("I116").Value = 6262.5713
$ws.Range("K116").Value = 6262.5713
$ws.Range("M116").Value = -2820.5713
$ws.Range("H132").Value = 265379.38
$ws.Range("I132").Value = 2195.9375
$ws.Range("K132").Value = 6587.8125
$ws.Range("M132").Value = -4057.8125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1475.2222
$ws.Range("I74").Value = 798.1429000000001
$ws.Range("J74").Value = 2204.3845
$ws.Range("K74").Value = 798.1429000000001
$ws.Range("L74").Value = 2204.3845
$ws.Range("M74").Value = 75.85709999999995
$ws.Range("N74").Value = -3952.3845
$ws.Range("H77").Value = 1475.2222
$ws.Range("I77").Value = 798.1429000000001
$ws.Range("J77").Value = 2204.3845
$ws.Range("K77").Value = 3990.7145
$ws.Range("L77").Value = 11021.9225
$ws.Range("M77").Value = 377.2855
$ws.Range("N77").Value = -19757.9225
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H122").Value = 2174.25
$ws.Range("I122").Value = 2199.8333
$ws.Range("K122").Value = 6599.499899999999
$ws.Range("M122").Value = -4149.499899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1852.8
$ws.Range("I107").Value = 1029.5
$ws.Range("J107").Value = 2401.6667
$ws.Range("K107").Value = 1029.5
$ws.Range("L107").Value = 2401.6667
$ws.Range("M107").Value = 890.5
$ws.Range("N107").Value = -6241.6667
$ws.Range("H134").Value = 3002.4614
$ws.Range("I134").Value = 2484.875
$ws.Range("K134").Value = 7454.625
$ws.Range("M134").Value = -4919.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3182.7778
$ws.Range("I31").Value = 2738.75
$ws.Range("J31").Value = 3538
$ws.Range("K31").Value = 2738.75
$ws.Range("L31").Value = 3538
$ws.Range("M31").Value = -2443.75
$ws.Range("N31").Value = -4128
$ws.Range("H34").Value = 3182.7778
$ws.Range("I34").Value = 2738.75
$ws.Range("J34").Value = 3538
$ws.Range("K34").Value = 2738.75
$ws.Range("L34").Value = 3538
$ws.Range("M34").Value = -2536.75
$ws.Range("N34").Value = -3942
$ws.Range("H99").Value = 36837.816
$ws.Range("I99").Value = 7246.3335
$ws.Range("J99").Value = 169999.5
$ws.Range("K99").Value = 7246.3335
$ws.Range("L99").Value = 169999.5
$ws.Range("M99").Value = -5748.3335
$ws.Range("N99").Value = -172995.5
$ws.Range("H107").Value = 1358.0834
$ws.Range("I107").Value = 1141.7333
$ws.Range("K107").Value = 1141.7333
$ws.Range("M107").Value = 778.2666999999999
$ws.Range("H126").Value = 36837.816
$ws.Range("I126").Value = 7246.3335
$ws.Range("J126").Value = 169999.5
$ws.Range("K126").Value = 21739.0005
$ws.Range("L126").Value = 509998.5
$ws.Range("M126").Value = -19269.0005
$ws.Range("N126").Value = -514938.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 70757160
$ws.Range("I4").Value = 4772241.5
$ws.Range("J4").Value = 466666700
$ws.Range("K4").Value = 14316724.5
$ws.Range("L4").Value = 1400000100
$ws.Range("M4").Value = -14316612.5
$ws.Range("N4").Value = -1400000324
$ws.Range("H5").Value = 1540.2354
$ws.Range("I5").Value = 1485.6
$ws.Range("J5").Value = 1618.2858
$ws.Range("K5").Value = 4456.799999999999
$ws.Range("L5").Value = 4854.857400000001
$ws.Range("M5").Value = -4344.799999999999
$ws.Range("N5").Value = -5078.857400000001
$ws.Range("H51").Value = 2256.75
$ws.Range("I51").Value = 2274
$ws.Range("J51").Value = 2251
$ws.Range("K51").Value = 6822
$ws.Range("L51").Value = 6753
$ws.Range("M51").Value = -6362
$ws.Range("N51").Value = -7673
$ws.Range("H60").Value = 397.16666
$ws.Range("I60").Value = 328.33334
$ws.Range("J60").Value = 466
$ws.Range("K60").Value = 985.0000200000001
$ws.Range("L60").Value = 1398
$ws.Range("M60").Value = -734.0000200000001
$ws.Range("N60").Value = -1900
$ws.Range("H76").Value = 3874.5
$ws.Range("I76").Value = 499
$ws.Range("J76").Value = 7250
$ws.Range("K76").Value = 1497
$ws.Range("L76").Value = 21750
$ws.Range("M76").Value = -1114
$ws.Range("N76").Value = -22516
$ws.Range("H79").Value = 3874.5
$ws.Range("I79").Value = 499
$ws.Range("J79").Value = 7250
$ws.Range("K79").Value = 1497
$ws.Range("L79").Value = 21750
$ws.Range("M79").Value = -171
$ws.Range("N79").Value = -24402
$ws.Range("H135").Value = 1540.2354
$ws.Range("I135").Value = 1485.6
$ws.Range("J135").Value = 1618.2858
$ws.Range("K135").Value = 13370.4
$ws.Range("L135").Value = 14564.5722
$ws.Range("M135").Value = -10835.4
$ws.Range("N135").Value = -19634.5722

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 11361.4
$ws.Range("I80").Value = 3370.7144
$ws.Range("J80").Value = 30006.334
$ws.Range("K80").Value = 3370.7144
$ws.Range("L80").Value = 30006.334
$ws.Range("M80").Value = -2372.7144
$ws.Range("N80").Value = -32002.334
$ws.Range("H83").Value = 11361.4
$ws.Range("I83").Value = 3370.7144
$ws.Range("J83").Value = 30006.334
$ws.Range("K83").Value = 16853.572
$ws.Range("L83").Value = 150031.67
$ws.Range("M83").Value = -11861.572
$ws.Range("N83").Value = -160015.67
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()
$ws.Range("H113").Value = 2900.111
$ws.Range("I113").Value = 1867.5
$ws.Range("K113").Value = 1867.5
$ws.Range("M113").Value = 302.5
$ws.Range("H122").Value = 3195.375
$ws.Range("I122").Value = 3109.2307
$ws.Range("K122").Value = 9327.6921
$ws.Range("M122").Value = -6877.6921
$ws.Range("H126").Value = 3982.8
$ws.Range("J126").Value = 4499.4287
$ws.Range("L126").Value = 13498.2861
$ws.Range("N126").Value = -18438.2861
$ws.Range("H132").Value = 3955.8462
$ws.Range("I132").Value = 2476.625
$ws.Range("K132").Value = 7429.875
$ws.Range("M132").Value = -4899.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2362.4167
$ws.Range("I68").Value = 2213.5454
$ws.Range("J68").Value = 4000
$ws.Range("K68").Value = 2213.5454
$ws.Range("L68").Value = 4000
$ws.Range("M68").Value = -1464.5454
$ws.Range("N68").Value = -5498
$ws.Range("H71").Value = 2362.4167
$ws.Range("I71").Value = 2213.5454
$ws.Range("J71").Value = 4000
$ws.Range("K71").Value = 11067.727
$ws.Range("L71").Value = 20000
$ws.Range("M71").Value = -7323.726999999999
$ws.Range("N71").Value = -27488
$ws.Range("H122").Value = 5696.8
$ws.Range("I122").Value = 5078.5
$ws.Range("K122").Value = 15235.5
$ws.Range("M122").Value = -12785.5
$ws.Range("H136").Value = 2799.8
$ws.Range("I136").Value = 1619.6
$ws.Range("K136").Value = 4858.799999999999
$ws.Range("M136").Value = -2308.799999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1324.7894
$ws.Range("I136").Value = 997.63635
$ws.Range("K136").Value = 2992.90905
$ws.Range("M136").Value = -442.9090500000002
